# CRM-1761 Add customer mobile no in pending booking and panding spare on partner pannel
#
# Inserts a new "Customer Phone Number" column right after "Customer Name"
# (new column B), pushing every subsequent column one position to the
# right. The header goes in row 1, the corresponding template placeholder
# {spare:customer_mobile} goes in row 2 - mirroring the existing
# "Service Center Mobile No" / {spare:primary_contact_phone_1} pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B - this shifts Invoice Id, Booking Id, etc.
# one column to the right, and the new column inherits column A's
# formatting (matches how Excel performs a manual "Insert Column").
$ws.Columns("B").Insert()

# New column header + template placeholder.
$ws.Range("B1").Value = "Customer Phone Number"
$ws.Range("B2").Value = "{spare:customer_mobile}"

# Match the column width that column A already has, same as what
# happens visually when a column is inserted next to it.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Leave the selection on the newly-populated cell.
$ws.Range("B2").Select() | Out-Null
